$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.406.25'
$ws.Range('E2').Value = '  -0.33%  '
$ws.Range('D3').Value = '1.847.28'
$ws.Range('E3').Value = '  -0.12%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9984'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.28'
$ws.Range('E5').Value = '  -0.78%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.6328'
$ws.Range('E6').Value = '  +0.63%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9998'
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07559'
$ws.Range('E8').Value = '  +0.20%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2968'
$ws.Range('E9').Value = '  -0.07%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '24.61'
$ws.Range('E10').Value = '  +1.61%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07731'
$ws.Range('D12').Value = '1.841.24'
$ws.Range('E12').Value = '  -2.02%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.002'
$ws.Range('E13').Value = '  -0.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.6853'
$ws.Range('E14').Value = '  +0.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001005'
$ws.Range('E15').Value = '  +2.56%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '83.17'
$ws.Range('E16').Value = '  -0.70%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.186'
$ws.Range('E17').Value = '  -0.38%  '
$ws.Range('D18').Value = '29.422.12'
$ws.Range('E18').Value = '  -0.42%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '230.31'
$ws.Range('E19').Value = '  -1.55%  '
$ws.Range('E20').Value = '  -0.11%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9996'
$ws.Range('E21').Value = '  -0.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.586'
$ws.Range('E22').Value = '  -0.28%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.000'
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '157.04'
$ws.Range('E24').Value = '  +0.82%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1402'
$ws.Range('E25').Value = '  +1.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.386'
$ws.Range('E26').Value = '  -0.53%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.70'
$ws.Range('E27').Value = '  -0.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.467'
$ws.Range('E28').Value = '  -1.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.05736'
$ws.Range('E29').Value = '  -1.79%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.251'
$ws.Range('E30').Value = '  -2.41%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.133'
$ws.Range('E31').Value = '  +0.59%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.043'
$ws.Range('E32').Value = '  +0.06%  '
$ws.Range('E33').Value = '  -1.92%  '
$ws.Range('E34').Value = '  -1.24%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7177'
$ws.Range('E35').Value = '  +0.26%  '
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('D37').Value = '1.253.14'
$ws.Range('E37').Value = '  +1.39%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01817'
$ws.Range('E38').Value = '  +2.20%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.780'
$ws.Range('E39').Value = '  -0.71%  '
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9106'
$ws.Range('E40').Value = '  -0.14%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.212'
$ws.Range('E41').Value = '  +1.38%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.000'
$ws.Range('E42').Value = '  +0.06%  '
$ws.Range('D43').Value = '2.002.27'
$ws.Range('E43').Value = '  -1.73%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '101.87'
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '66.48'
$ws.Range('E45').Value = '  -1.43%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000119'
$ws.Range('E46').Value = '  +1.70%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.066'
$ws.Range('E47').Value = '  -2.96%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.160'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4036'
$ws.Range('E49').Value = '  +0.20%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.711'
$ws.Range('E50').Value = '  +0.68%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.1132'
$ws.Range('E51').Value = '  +1.03%  '
